# Sync attendance_reports: fix "Recorded By" (column G) ordering so the
# human reviewer/editor email is listed before the literal "System" tag
# (and similar reordering for the few multi-editor values).
#
# The workbook has a single worksheet ("Session Analysis Results") with a
# header row (row 1) and data rows 2..157. Column G = "Recorded By".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact value -> value remap observed between the two revisions. Using an
# explicit lookup (rather than a generic "reverse the list" transform)
# because one combination - "admin@admin.com, System" - is intentionally
# left untouched while every other combination is reordered.
$map = @{
    "System, backup@backdoor.com, system" = "system, backup@backdoor.com, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = 7
    $current = $cell.Value2

    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
